# BigAssTableOfLenses.xlsx edit
# Commit: "+CP.3, +GF 23mm, GF 110mm" -- restore the Zeiss CP.3 XD compact
# prime lineup and add the Fujifilm GF 23mm f/4 and GF 110mm f/2 lenses,
# mark a batch of previously-unverified "Designed In"/"Made In" cells with
# the red "Bad" highlight, and log the change in the Changelog sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LensTable")
$cl = $wb.Worksheets.Item("Changelog")

# ---------------------------------------------------------------------
# 1. Mark previously blank Designed-In/Made-In cells as "Bad" (missing
#    data) across the rows that didn't yet have them flagged.
# ---------------------------------------------------------------------

$pOnlyRows = @(89,124,125,126,127,128,129,130,131,132,133,134,135,136,137,138,139,140,141,
  147,148,149,150,151,152,153,154,155,156,157,158,159,160,161,162,163,164,165,166,
  182,183,184,185,186,187,188,189,190,191,192,193,194,195,196,197,198,199,200,201,
  202,203,204,205,206,207,208,209,210,211,212,213,214,215,216,217,218,219,220,221,
  222,223,224,225,226,227,228,229,230,231,232,
  317,318,319,320,
  364,365,366,367,368,369,370,371,372,373,374,375,376,377,378,
  385,387,388,389,
  400,401,402,403,404)

foreach ($r in $pOnlyRows) {
  $ws.Cells.Item($r, 16).Style = "Bad"
}

# Rows 142-146 already had a blank (plain-styled) "Made In" cell -- just
# recolor it to the "Bad" highlight to match the rest of the column.
foreach ($r in @(142,143,144,145,146)) {
  $ws.Cells.Item($r, 16).Style = "Bad"
}

# Rows that need both Designed-In (O) and Made-In (P) flagged as "Bad".
foreach ($r in @(90,91,396,397)) {
  $ws.Cells.Item($r, 15).Style = "Bad"
  $ws.Cells.Item($r, 16).Style = "Bad"
}

# Rows 237-253: same, but also boxed with a thin grey border.
foreach ($r in 237..253) {
  $rng = $ws.Range("O" + $r + ":P" + $r)
  $rng.Style = "Bad"
  $rng.Borders.LineStyle = 1
  $rng.Borders.Weight = 2
  $rng.Borders.Color = 8355711
}

# ---------------------------------------------------------------------
# 2. Restore the Zeiss CP.3 XD compact prime lineup (rows 407-415) and
#    add the two new Fujifilm GF lenses (rows 416-417).
# ---------------------------------------------------------------------

function Set-LensRow {
  param(
    [int]$Row,
    [string]$Manufacture,
    [string]$Model,
    [int]$ReleaseYear,
    [double]$MSRP,
    [double]$Weight,
    [double]$Length,
    [double]$MaxDiameter,
    $MaxAperture,
    [string]$EFL,
    [string]$Zoom,
    [int]$Stabilized,
    [int]$Cinema,
    [string]$Format,
    [string]$Variety,
    [string]$DesignedIn,
    [string]$MadeIn,
    [int]$Index
  )

  $ws.Cells.Item($Row, 1).Value = $Manufacture
  $ws.Cells.Item($Row, 2).Value = $Model
  $ws.Cells.Item($Row, 3).Value = $ReleaseYear
  $ws.Cells.Item($Row, 4).Value = $MSRP
  $ws.Cells.Item($Row, 5).Value = $Weight
  $ws.Cells.Item($Row, 6).Value = $Length
  $ws.Cells.Item($Row, 7).Value = $MaxDiameter
  $ws.Cells.Item($Row, 8).Value = $MaxAperture
  $ws.Cells.Item($Row, 9).NumberFormat = "@"
  $ws.Cells.Item($Row, 9).Value = $EFL
  $ws.Cells.Item($Row, 10).NumberFormat = "@"
  $ws.Cells.Item($Row, 10).Value = $Zoom
  $ws.Cells.Item($Row, 11).Value = $Stabilized
  $ws.Cells.Item($Row, 12).Value = $Cinema
  $ws.Cells.Item($Row, 13).Value = $Format
  $ws.Cells.Item($Row, 14).Value = $Variety
  $ws.Cells.Item($Row, 15).Value = $DesignedIn
  if ($MadeIn -ne $null) {
    $ws.Cells.Item($Row, 16).Value = $MadeIn
  }
  $ws.Cells.Item($Row, 17).Value = $Index
}

Set-LensRow -Row 407 -Manufacture "Zeiss" -Model "CP.3 XD 100mm T2.1 Compact Prime" `
  -ReleaseYear 2017 -MSRP 6690 -Weight 1010 -Length 126.5 -MaxDiameter 95 -MaxAperture 2 `
  -EFL "100" -Zoom "0" -Stabilized 0 -Cinema 1 -Format "Full-Frame" -Variety "Telephoto" `
  -DesignedIn "Germany" -MadeIn "Germany" -Index 406

Set-LensRow -Row 408 -Manufacture "Zeiss" -Model "CP.3 XD 15mm T2.9 Compact Prime" `
  -ReleaseYear 2017 -MSRP 7490 -Weight 870 -Length 83.7 -MaxDiameter 95 -MaxAperture 2.8 `
  -EFL "15" -Zoom "0" -Stabilized 0 -Cinema 1 -Format "Full-Frame" -Variety "Ultra Wide Angle" `
  -DesignedIn "Germany" -MadeIn "Germany" -Index 407

Set-LensRow -Row 409 -Manufacture "Zeiss" -Model "CP.3 XD 18mm T2.9 Compact Prime" `
  -ReleaseYear 2017 -MSRP 6690 -Weight 860 -Length 83.7 -MaxDiameter 95 -MaxAperture 2.8 `
  -EFL "18" -Zoom "0" -Stabilized 0 -Cinema 1 -Format "Full-Frame" -Variety "Ultra Wide Angle" `
  -DesignedIn "Germany" -MadeIn "Germany" -Index 408

Set-LensRow -Row 410 -Manufacture "Zeiss" -Model "CP.3 XD 21mm T2.9 Compact Prime" `
  -ReleaseYear 2017 -MSRP 5790 -Weight 820 -Length 83.7 -MaxDiameter 95 -MaxAperture 2.8 `
  -EFL "21" -Zoom "0" -Stabilized 0 -Cinema 1 -Format "Full-Frame" -Variety "Ultra Wide Angle" `
  -DesignedIn "Germany" -MadeIn "Germany" -Index 409

Set-LensRow -Row 411 -Manufacture "Zeiss" -Model "CP.3 XD 25mm T2.1 Compact Prime" `
  -ReleaseYear 2017 -MSRP 5790 -Weight 820 -Length 83.7 -MaxDiameter 95 -MaxAperture 2 `
  -EFL "25" -Zoom "0" -Stabilized 0 -Cinema 1 -Format "Full-Frame" -Variety "Wide Angle" `
  -DesignedIn "Germany" -MadeIn "Germany" -Index 410

Set-LensRow -Row 412 -Manufacture "Zeiss" -Model "CP.3 XD 28mm T2.1 Compact Prime" `
  -ReleaseYear 2017 -MSRP 5790 -Weight 840 -Length 83.7 -MaxDiameter 95 -MaxAperture 2 `
  -EFL "28" -Zoom "0" -Stabilized 0 -Cinema 1 -Format "Full-Frame" -Variety "Wide Angle" `
  -DesignedIn "Germany" -MadeIn "Germany" -Index 411

Set-LensRow -Row 413 -Manufacture "Zeiss" -Model "CP.3 XD 35mm T2.1 Compact Prime" `
  -ReleaseYear 2017 -MSRP 5790 -Weight 800 -Length 83.7 -MaxDiameter 95 -MaxAperture 2 `
  -EFL "35" -Zoom "0" -Stabilized 0 -Cinema 1 -Format "Full-Frame" -Variety "Wide Angle" `
  -DesignedIn "Germany" -MadeIn "Germany" -Index 412

Set-LensRow -Row 414 -Manufacture "Zeiss" -Model "CP.3 XD 50mm T2.1 Compact Prime" `
  -ReleaseYear 2017 -MSRP 5790 -Weight 770 -Length 83.7 -MaxDiameter 95 -MaxAperture 2 `
  -EFL "50" -Zoom "0" -Stabilized 0 -Cinema 1 -Format "Full-Frame" -Variety "Normal" `
  -DesignedIn "Germany" -MadeIn "Germany" -Index 413

Set-LensRow -Row 415 -Manufacture "Zeiss" -Model "CP.3 XD 85mm T2.1 Compact Prime" `
  -ReleaseYear 2017 -MSRP 5790 -Weight 880 -Length 83.7 -MaxDiameter 95 -MaxAperture 2 `
  -EFL "85" -Zoom "0" -Stabilized 0 -Cinema 1 -Format "Full-Frame" -Variety "Short Telephoto" `
  -DesignedIn "Germany" -MadeIn "Germany" -Index 414

Set-LensRow -Row 416 -Manufacture "Fujifilm" -Model "GF 23mm f/4 R LM WR" `
  -ReleaseYear 2017 -MSRP 2599 -Weight 845 -Length 103 -MaxDiameter 89.8 -MaxAperture 4 `
  -EFL "23" -Zoom "0" -Stabilized 0 -Cinema 0 -Format "Crop Medium Format Digital" -Variety "Ultra Wide Angle" `
  -DesignedIn "Japan" -MadeIn $null -Index 415
$ws.Cells.Item(416, 16).Style = "Bad"

Set-LensRow -Row 417 -Manufacture "Fujifilm" -Model "GF 110mm f/2 R LM WR" `
  -ReleaseYear 2017 -MSRP 2799 -Weight 1010 -Length 125.5 -MaxDiameter 94.3 -MaxAperture 2 `
  -EFL "110" -Zoom "0" -Stabilized 0 -Cinema 0 -Format "Crop Medium Format Digital" -Variety "Short Telephoto" `
  -DesignedIn "Japan" -MadeIn $null -Index 416

# ---------------------------------------------------------------------
# 3. Extend Table1 to cover the new rows.
# ---------------------------------------------------------------------

$ws.ListObjects.Item("Table1").Resize($ws.Range("A1:Q417"))

# ---------------------------------------------------------------------
# 4. Log the change in the Changelog sheet.
# ---------------------------------------------------------------------

$cl.Cells.Item(53, 1).NumberFormat = "m/d/yyyy"
$cl.Cells.Item(53, 1).Value = "6/28/2017"
$cl.Cells.Item(53, 2).Value = "Restore CP.3 lenses, + GF 23mm and GF 110mm"

# ---------------------------------------------------------------------
# 5. Selection / active sheet bookkeeping (LensTable becomes the active
#    tab, matching where editing left off).
# ---------------------------------------------------------------------

$cl.Range("B54").Select()
$ws.Activate()
$ws.Range("E417").Select()
